$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing account 005046790 / BEATRIZ (Excel row 3)
$ws.Rows.Item(3).Delete()
